$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header row: "_old" columns -> "_FV2410", "_new" columns -> "_FV2504" ---
$oldPrefixes = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

# Columns A:J -> "<name>_FV2410"
for ($i = 0; $i -lt $oldPrefixes.Length; $i++) {
    $col = [char](65 + $i)  # A..J
    $ws.Range("$col" + "1").Value = $oldPrefixes[$i] + "_FV2410"
}

# Column K stays "diff" (unchanged)

# Columns L:U -> "<name>_FV2504"
for ($i = 0; $i -lt $oldPrefixes.Length; $i++) {
    $col = [char](76 + $i)  # L..U
    $ws.Range("$col" + "1").Value = $oldPrefixes[$i] + "_FV2504"
}

# --- Add Table1 over A1:U62 with autofilter, banded rows, no explicit style name ---
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:U62"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
$tbl.ShowTableStyleFirstColumn = $false
$tbl.ShowTableStyleLastColumn = $false
$tbl.ShowTableStyleRowStripes = $true
$tbl.ShowTableStyleColumnStripes = $false
$tbl.TableStyle = ""

# --- Freeze header row (pane split after row 1) ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
